# Updated additional twitter handle
# Locate the subtitle shape on the "Scheduling Tweets with PowerShell" title
# slide that holds the contact-info paragraphs (identified by the
# "Github.com/ergo3114" text, which is unique in the deck).
$p = $ppt.ActivePresentation

$targetSlide = $null
$targetShape = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $sl = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $sl.Shapes.Count; $shi++) {
        $shp = $sl.Shapes.Item($shi)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -like "*ergo3114*") {
                $targetSlide = $sl
                $targetShape = $shp
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange

# Find the paragraph that reads "Github.com/ergo3114" - the new "Twitter:
# @PSHChatt" paragraph needs to be inserted immediately before it.
$paraCount = $tr.Paragraphs().Count
$githubParaIndex = -1
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $tr.Paragraphs($i, 1)
    if ($para.Text -eq "Github.com/ergo3114") {
        $githubParaIndex = $i
    }
}

$githubPara = $tr.Paragraphs($githubParaIndex, 1)

# Insert a brand-new paragraph right before it containing the new handle.
$newPara = $githubPara.InsertBefore("Twitter: @PSHChatt" + [char]13)

# Re-fetch the freshly inserted paragraph so further edits address it
# directly (InsertBefore's return range can shift after the text mutates).
$newPara = $tr.Paragraphs($githubParaIndex, 1)

# Split "Twitter: @PSHChatt" into two runs - "Twitter: @" and "PSHChatt" -
# matching the two-run structure of the existing "Twitter: @littlejohnpsh"
# paragraph above it. Re-stating the (unchanged) font name on the back half
# of the text forces the engine to materialize it as its own run.
$handleChars = $newPara.Characters(11, 8)
$handleChars.Font.Name = $handleChars.Font.Name
